$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("H1").Value = "store_count"

# Row 2
$ws.Range("E2").Value = "80170100"
$ws.Range("F2").Value = "800 E Dimond Blvd"
$ws.Range("G2").Value = 9075220666
$ws.Range("H2").Value = 50

# Row 3
$ws.Range("E3").Value = "3138"
$ws.Range("F3").Value = "5725 Johnston St."
$ws.Range("G3").Value = "337-993-1090"
$ws.Range("H3").Value = 20

# Row 4
$ws.Range("E4").Value = "68846"
$ws.Range("F4").Value = "6470 SPALDING DR"
$ws.Range("G4").Value = "(770) 582-0841"
$ws.Range("H4").Value = 20

# Row 5
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "5211 E. Washington Blvd."
$ws.Range("G5").Value = "562-463-9222"
$ws.Range("H5").Value = 20

# Row 6
$ws.Range("E6").Value = "401"
$ws.Range("F6").Value = "3200 Rogers Ave."
$ws.Range("G6").Value = "(479) 709-8800"
$ws.Range("H6").Value = 20

# selection change
$ws.Range("C9").Select()
